$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A6 with the combined tuple-style strings
$ws.Range("A2").Value = "('森', ['Basic Land — Forest', '({T}: Add {G}.)'])"
$ws.Range("A3").Value = "('島', ['Basic Land — Island', '({T}: Add {U}.)'])"
$ws.Range("A4").Value = "('山', ['Basic Land — Mountain', '({T}: Add {R}.)'])"
$ws.Range("A5").Value = "('平地', ['Basic Land — Plains', '({T}: Add {W}.)'])"
$ws.Range("A6").Value = "('沼', ['Basic Land — Swamp', '({T}: Add {B}.)'])"

# Remove now-unused rows 7-16
$ws.Range("A7:A16").EntireRow.Delete()
